$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value for D (Price) column; written as text to avoid numeric auto-conversion
$dValues = @{
    'D2' = '27.278.15'
    'D3' = '1.851.04'
    'D4' = '1.002'
    'D5' = '324.25'
    'D7' = '0.4533'
    'D8' = '0.3877'
    'D9' = '48.49'
    'D10' = '0.07909'
    'D11' = '1.012'
    'D13' = '1.833.46'
    'D14' = '5.891'
    'D15' = '7.129'
    'D17' = '0.06591'
    'D18' = '85.70'
    'D19' = '0.00001022'
    'D20' = '17.12'
    'D22' = '5.481'
    'D23' = '27.290.52'
    'D24' = '10.86'
    'D25' = '2.302'
    'D26' = '2.071.85'
    'D27' = '153.93'
    'D28' = '19.90'
    'D29' = '2.048'
    'D30' = '5.414'
    'D31' = '121.18'
    'D32' = '0.09297'
    'D33' = '0.9357'
    'D34' = '1.443'
    'D35' = '3.587'
    'D36' = '5.241'
    'D37' = '0.06027'
    'D38' = '0.02220'
    'D39' = '1.205'
    'D40' = '8.078'
    'D42' = '0.5899'
    'D43' = '0.1876'
    'D45' = '1.266'
    'D46' = '0.5591'
    'D47' = '12.08'
    'D49' = '1.905'
    'D50' = '0.06730'
    'D51' = '107.76'
}

# Map of cell -> new value for E (Volume 1h) column
$eValues = @{
    'E2' = '  -3.14%  '
    'E3' = '  -3.96%  '
    'E4' = '  +0.18%  '
    'E5' = '  -1.93%  '
    'E6' = '  +0.08%  '
    'E7' = '  -4.08%  '
    'E8' = '  -4.52%  '
    'E9' = '  -8.62%  '
    'E10' = '  -6.25%  '
    'E11' = '  -3.63%  '
    'E12' = '  -4.50%  '
    'E13' = '  -5.15%  '
    'E14' = '  -3.46%  '
    'E15' = '  -5.13%  '
    'E16' = '  +0.05%  '
    'E17' = '  -0.06%  '
    'E18' = '  -5.40%  '
    'E19' = '  -3.94%  '
    'E20' = '  -5.68%  '
    'E21' = '  -0.02%  '
    'E22' = '  -4.69%  '
    'E23' = '  -3.06%  '
    'E24' = '  -4.69%  '
    'E25' = '  +0.78%  '
    'E26' = '  -4.08%  '
    'E27' = '  -0.14%  '
    'E28' = '  -1.02%  '
    'E29' = '  -4.93%  '
    'E30' = '  -5.67%  '
    'E31' = '  -2.07%  '
    'E32' = '  -3.30%  '
    'E33' = '  -4.41%  '
    'E34' = '  -0.98%  '
    'E35' = '  -1.28%  '
    'E36' = '  -5.96%  '
    'E37' = '  -2.53%  '
    'E38' = '  -4.30%  '
    'E39' = '  -2.48%  '
    'E40' = '  -10.53%  '
    'E41' = '  +0.07%  '
    'E42' = '  -4.44%  '
    'E43' = '  -1.56%  '
    'E44' = '  -8.83%  '
    'E45' = '  -3.00%  '
    'E46' = '  -5.09%  '
    'E47' = '  -5.61%  '
    'E48' = '  -2.90%  '
    'E49' = '  -6.47%  '
    'E50' = '  -1.57%  '
    'E51' = '  -2.25%  '
}

foreach ($cellRef in $dValues.Keys) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $dValues[$cellRef]
    $cell.Style = $origStyle
}

foreach ($cellRef in $eValues.Keys) {
    $ws.Range($cellRef).Value = $eValues[$cellRef]
}
